$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1504
$ws.Range("F5").Value = 7640
$ws.Range("F6").Value = 4833
$ws.Range("F7").Value = 7117
$ws.Range("F9").Value = 286
$ws.Range("F10").Value = 1510
$ws.Range("F11").Value = 870
$ws.Range("F12").Value = 200
$ws.Range("F13").Value = 63
$ws.Range("F14").Value = 1173
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 8
$ws.Range("F20").Value = 236
$ws.Range("F22").Value = 1195
$ws.Range("F25").Value = 48
$ws.Range("F26").Value = 1242
$ws.Range("F27").Value = 49
$ws.Range("F28").Value = 154
$ws.Range("F31").Value = 202
$ws.Range("F32").Value = 21
$ws.Range("F33").Value = 5
$ws.Range("F34").Value = 50
$ws.Range("F35").Value = 116
$ws.Range("F37").Value = 559
$ws.Range("F38").Value = 429
$ws.Range("F39").Value = 79
$ws.Range("F41").Value = 93
$ws.Range("F42").Value = 402
$ws.Range("F44").Value = 588
$ws.Range("F45").Value = 148

$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 134
$ws.Range("F25").Value = 7
$ws.Range("F27").Value = 640
$ws.Range("F29").Value = 32
$ws.Range("F32").Value = 867
$ws.Range("F34").Value = 611
$ws.Range("F37").Value = 115
$ws.Range("F41").Value = 146
$ws.Range("F47").Value = 9

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 734
$ws.Range("F6").Value = 685
$ws.Range("F8").Value = 61
$ws.Range("F9").Value = 68
$ws.Range("F10").Value = 1662
$ws.Range("F11").Value = 2568

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 734
$ws.Range("F3").Value = 1504
$ws.Range("F7").Value = 685
$ws.Range("F8").Value = 685
$ws.Range("F9").Value = 7640
$ws.Range("F10").Value = 4834
$ws.Range("F11").Value = 7117
$ws.Range("F12").Value = 286
$ws.Range("F13").Value = 1510
$ws.Range("F15").Value = 870
$ws.Range("F17").Value = 200
$ws.Range("F18").Value = 1662
$ws.Range("F19").Value = 2568
$ws.Range("F21").Value = 63
$ws.Range("F22").Value = 1173
$ws.Range("F25").Value = 236
$ws.Range("F26").Value = 1195
$ws.Range("F27").Value = 640
$ws.Range("F29").Value = 1242
$ws.Range("F30").Value = 154
$ws.Range("F31").Value = 202
$ws.Range("F32").Value = 32
$ws.Range("F33").Value = 867
$ws.Range("F34").Value = 50
$ws.Range("F35").Value = 116
$ws.Range("F37").Value = 559
$ws.Range("F38").Value = 611
$ws.Range("F39").Value = 79
$ws.Range("F41").Value = 93
$ws.Range("F42").Value = 115
$ws.Range("F43").Value = 402
$ws.Range("F44").Value = 588
$ws.Range("F47").Value = 148
